# Add 2022-Q4 data:
#   1) Insert a new worksheet "2022-Q4" right before "2022-Q3" (after "总计").
#   2) Populate the new sheet with the Q4 fund-holding table.
#   3) Insert a new row into the "总计" (summary) sheet for the 2022-Q4 period,
#      shifting the existing rows down and re-numbering the index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert the new "2022-Q4" sheet before "2022-Q3"
# ---------------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ4 = $wb.Worksheets.Add($wsQ3)
$wsQ4.Name = "2022-Q4"

# ---------------------------------------------------------------------------
# 2) Populate "2022-Q4" sheet
# ---------------------------------------------------------------------------
$wsQ4.Range("B1").Value = "基金代码"
$wsQ4.Range("C1").Value = "基金名称"
$wsQ4.Range("D1").Value = "基金规模"
$wsQ4.Range("E1").Value = "股票总仓位"
$wsQ4.Range("F1").Value = "仓位占比"
$wsQ4.Range("G1").Value = "持有市值(亿元)"
$wsQ4.Range("H1").Value = "仓位排名"

foreach ($cell in @("B1","C1","D1","E1","F1","G1","H1")) {
    $r = $wsQ4.Range($cell)
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
    $r.Borders.LineStyle = 1
}

$q4Rows = @(
    @("000593", "易方达标普全球高端消费品指数增强（QDII）美元现汇", "2.30", "93.71", "6.13", "0.1410", 7),
    @("005676", "易方达标普全球高端消费品指数增强C（QDII）人民币",   "2.30", "93.71", "6.13", "0.1410", 7),
    @("118002", "易方达标普全球高端消费品指数增强A（QDII）人民币",   "2.30", "93.71", "6.13", "0.1410", 7)
)

$r = 2
foreach ($row in $q4Rows) {
    # Leading "'" keeps numeric-looking text (fund codes, 2-decimal figures)
    # stored as text instead of Excel auto-coercing it to a number.
    $wsQ4.Cells.Item($r, 1).Value = $r - 2
    $wsQ4.Cells.Item($r, 2).Value = "'" + $row[0]
    $wsQ4.Cells.Item($r, 3).Value = $row[1]
    $wsQ4.Cells.Item($r, 4).Value = "'" + $row[2]
    $wsQ4.Cells.Item($r, 5).Value = "'" + $row[3]
    $wsQ4.Cells.Item($r, 6).Value = "'" + $row[4]
    $wsQ4.Cells.Item($r, 7).Value = "'" + $row[5]
    $wsQ4.Cells.Item($r, 8).Value = $row[6]
    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Insert the 2022-Q4 row into the "总计" summary sheet
# ---------------------------------------------------------------------------
$wsSum = $wb.Worksheets.Item("总计")
$wsSum.Rows.Item(2).Insert()
$wsSum.Range("B2:D2").ClearFormats()

$a2 = $wsSum.Cells.Item(2, 1)
$a2.Value = 0
$a2.Font.Bold = $true
$a2.HorizontalAlignment = -4108
$a2.VerticalAlignment = -4160
$a2.Borders.LineStyle = 1

$wsSum.Cells.Item(2, 2).Value = "2022-Q4"
$wsSum.Cells.Item(2, 3).Value = 3
$wsSum.Cells.Item(2, 4).Value = 0.42

# Re-number the index column (A) for the rows that shifted down one position.
for ($row = 3; $row -le 7; $row++) {
    $wsSum.Cells.Item($row, 1).Value = $row - 2
}
